# Apply the Mon May  8 11:07:16 UTC 2023 "cryptos" price/volume refresh.
# Updates columns D (Price) and E (Volume(1h)) text values for the changed
# rows, plus the RenderToken/PaxDollar row swap at rows 50-51 (B, C, D, E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.962.88'
$ws.Range("E2").Value = '  -3.29%  '
# Row 3
$ws.Range("D3").Value = '1.858.62'
$ws.Range("E3").Value = '  -2.68%  '
# Row 4
$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  -0.08%  '
# Row 5
$ws.Range("D5").Value = '''318.03'
$ws.Range("E5").Value = '  -2.14%  '
# Row 6
$ws.Range("D6").Value = '''1.000'
$ws.Range("E6").Value = '  -0.07%  '
# Row 7
$ws.Range("D7").Value = '''0.4371'
$ws.Range("E7").Value = '  -4.72%  '
# Row 8
$ws.Range("D8").Value = '''0.3700'
$ws.Range("E8").Value = '  -3.05%  '
# Row 9
$ws.Range("D9").Value = '''0.07496'
# Row 10
$ws.Range("D10").Value = '''0.9383'
$ws.Range("E10").Value = '  -4.32%  '
# Row 11
$ws.Range("D11").Value = '''21.25'
$ws.Range("E11").Value = '  -4.44%  '
# Row 12
$ws.Range("D12").Value = '1.861.01'
$ws.Range("E12").Value = '  -0.56%  '
# Row 13
$ws.Range("D13").Value = '''6.721'
$ws.Range("E13").Value = '  -3.40%  '
# Row 14
$ws.Range("D14").Value = '''5.439'
$ws.Range("E14").Value = '  -4.34%  '
# Row 15
$ws.Range("D15").Value = '''0.06845'
$ws.Range("E15").Value = '  -3.11%  '
# Row 16
$ws.Range("D16").Value = '''1.002'
$ws.Range("E16").Value = '  -0.09%  '
# Row 17
$ws.Range("D17").Value = '''81.59'
$ws.Range("E17").Value = '  -2.70%  '
# Row 18
$ws.Range("D18").Value = '''0.000009043'
$ws.Range("E18").Value = '  -4.22%  '
# Row 20
$ws.Range("D20").Value = '''15.94'
$ws.Range("E20").Value = '  -4.09%  '
# Row 21
$ws.Range("D21").Value = '27.954.45'
$ws.Range("E21").Value = '  -3.32%  '
# Row 22
$ws.Range("D22").Value = '''5.112'
$ws.Range("E22").Value = '  -3.96%  '
# Row 23
$ws.Range("E23").Value = '  +1.33%  '
# Row 24
$ws.Range("D24").Value = '2.062.83'
$ws.Range("E24").Value = '  -1.77%  '
# Row 25
$ws.Range("E25").Value = '  -4.45%  '
# Row 26
$ws.Range("D26").Value = '''154.74'
$ws.Range("E26").Value = '  -2.51%  '
# Row 27
$ws.Range("E27").Value = '  -3.54%  '
# Row 28
$ws.Range("D28").Value = '''5.434'
$ws.Range("E28").Value = '  -4.41%  '
# Row 29
$ws.Range("D29").Value = '''113.58'
$ws.Range("E29").Value = '  -3.37%  '
# Row 30
$ws.Range("D30").Value = '''1.738'
$ws.Range("E30").Value = '  -7.09%  '
# Row 31
$ws.Range("D31").Value = '''0.08982'
$ws.Range("E31").Value = '  -3.38%  '
# Row 32
$ws.Range("D32").Value = '''0.8117'
$ws.Range("E32").Value = '  -6.37%  '
# Row 33
$ws.Range("E33").Value = '  -5.62%  '
# Row 34
$ws.Range("D34").Value = '''1.175'
# Row 35
$ws.Range("D35").Value = '''2.936'
$ws.Range("E35").Value = '  -2.59%  '
# Row 36
$ws.Range("E36").Value = '  -0.04%  '
# Row 37
$ws.Range("D37").Value = '''0.05492'
$ws.Range("E37").Value = '  -3.98%  '
# Row 38
$ws.Range("D38").Value = '''1.109'
$ws.Range("E38").Value = '  -3.90%  '
# Row 39
$ws.Range("E39").Value = '  -3.48%  '
# Row 40
$ws.Range("D40").Value = '''2.914'
$ws.Range("E40").Value = '  +2.11%  '
# Row 41
$ws.Range("D41").Value = '''0.5259'
$ws.Range("E41").Value = '  -4.36%  '
# Row 42
$ws.Range("D42").Value = '''7.026'
$ws.Range("E42").Value = '  -5.56%  '
# Row 43
$ws.Range("D43").Value = '''0.1687'
$ws.Range("E43").Value = '  -3.73%  '
# Row 44
$ws.Range("D44").Value = '''8.808'
$ws.Range("E44").Value = '  -5.82%  '
# Row 45
$ws.Range("D45").Value = '''0.06788'
$ws.Range("E45").Value = '  -1.73%  '
# Row 46
$ws.Range("D46").Value = '''0.4892'
$ws.Range("E46").Value = '  -5.48%  '
# Row 47
$ws.Range("D47").Value = '''10.68'
$ws.Range("E47").Value = '  -4.94%  '
# Row 48
$ws.Range("D48").Value = '''106.43'
# Row 49
$ws.Range("D49").Value = '''1.681'
$ws.Range("E49").Value = '  -5.45%  '
# Row 50
$ws.Range("B50").Value = 'PaxDollar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D50").Value = '''0.9997'
$ws.Range("E50").Value = '  -0.15%  '
# Row 51
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").Value = '''1.906'
$ws.Range("E51").Value = '  -11.84%  '
